# Auto-generated: refresh market-price derived columns (H:N) for Leve profit sheets.
# Mirrors a scheduled data-refresh commit that updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns with newly polled Market Board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 11555.407
$ws.Range("I62").Value = 15993.066
$ws.Range("J62").Value = 6008.3335
$ws.Range("K62").Value = 15993.066
$ws.Range("L62").Value = 6008.3335
$ws.Range("M62").Value = -15369.066
$ws.Range("N62").Value = -7256.3335
$ws.Range("H65").Value = 11555.407
$ws.Range("I65").Value = 15993.066
$ws.Range("J65").Value = 6008.3335
$ws.Range("K65").Value = 79965.33
$ws.Range("L65").Value = 30041.6675
$ws.Range("M65").Value = -76845.33
$ws.Range("N65").Value = -36281.6675
$ws.Range("H116").Value = 3922.75
$ws.Range("J116").Value = 3513.1667
$ws.Range("L116").Value = 3513.1667
$ws.Range("N116").Value = -10397.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 772
$ws.Range("I2").Value = 670.3333
$ws.Range("K2").Value = 670.3333
$ws.Range("M2").Value = -557.3333
$ws.Range("H45").Value = 851.3077
$ws.Range("I45").Value = 763
$ws.Range("J45").Value = 1050
$ws.Range("K45").Value = 763
$ws.Range("L45").Value = 1050
$ws.Range("M45").Value = -386
$ws.Range("N45").Value = -1804
$ws.Range("H74").Value = 12787576
$ws.Range("I74").Value = 10098787
$ws.Range("J74").Value = 17599092
$ws.Range("K74").Value = 10098787
$ws.Range("L74").Value = 17599092
$ws.Range("M74").Value = -10097913
$ws.Range("N74").Value = -17600840
$ws.Range("H77").Value = 12787576
$ws.Range("I77").Value = 10098787
$ws.Range("J77").Value = 17599092
$ws.Range("K77").Value = 50493935
$ws.Range("L77").Value = 87995460
$ws.Range("M77").Value = -50489567
$ws.Range("N77").Value = -88004196
$ws.Range("H116").Value = 772
$ws.Range("I116").Value = 670.3333
$ws.Range("K116").Value = 670.3333
$ws.Range("M116").Value = 1623.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 772
$ws.Range("I3").Value = 670.3333
$ws.Range("K3").Value = 670.3333
$ws.Range("M3").Value = -556.3333
$ws.Range("H107").Value = 945.7273
$ws.Range("I107").Value = 945.7273
$ws.Range("K107").Value = 945.7273
$ws.Range("M107").Value = 974.2727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 929599.75
$ws.Range("I31").Value = 747.6724
$ws.Range("J31").Value = 2276435.2
$ws.Range("K31").Value = 747.6724
$ws.Range("L31").Value = 2276435.2
$ws.Range("M31").Value = -452.6724
$ws.Range("N31").Value = -2277025.2
$ws.Range("H34").Value = 929599.75
$ws.Range("I34").Value = 747.6724
$ws.Range("J34").Value = 2276435.2
$ws.Range("K34").Value = 747.6724
$ws.Range("L34").Value = 2276435.2
$ws.Range("M34").Value = -545.6724
$ws.Range("N34").Value = -2276839.2
$ws.Range("H99").Value = 63937
$ws.Range("I99").Value = 63937
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 63937
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -62439
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 1182.5483
$ws.Range("I105").Value = 1005.95
$ws.Range("J105").Value = 1503.6364
$ws.Range("K105").Value = 1005.95
$ws.Range("L105").Value = 1503.6364
$ws.Range("M105").Value = 741.05
$ws.Range("N105").Value = -4997.6364
$ws.Range("H107").Value = 290.77777
$ws.Range("I107").Value = 262.9
$ws.Range("J107").Value = 325.625
$ws.Range("K107").Value = 262.9
$ws.Range("L107").Value = 325.625
$ws.Range("M107").Value = 1657.1
$ws.Range("N107").Value = -4165.625
$ws.Range("H126").Value = 63937
$ws.Range("I126").Value = 63937
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 191811
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -189341
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 17859904
$ws.Range("I132").Value = 33335268
$ws.Range("J132").Value = 3716.2307
$ws.Range("K132").Value = 100005804
$ws.Range("L132").Value = 11148.6921
$ws.Range("M132").Value = -100003274
$ws.Range("N132").Value = -16208.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 214.58333
$ws.Range("I33").Value = 171.06667
$ws.Range("J33").Value = 287.1111
$ws.Range("K33").Value = 1026.40002
$ws.Range("L33").Value = 1722.6666
$ws.Range("M33").Value = -743.40002
$ws.Range("N33").Value = -2288.6666
$ws.Range("H47").Value = 240
$ws.Range("I47").Value = 240
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 720
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -289
$ws.Range("N47").ClearContents()
$ws.Range("H48").Value = 2580
$ws.Range("I48").Value = 500
$ws.Range("J48").Value = 3471.4285
$ws.Range("K48").Value = 1500
$ws.Range("L48").Value = 10414.2855
$ws.Range("M48").Value = -1250
$ws.Range("N48").Value = -10914.2855
$ws.Range("H68").Value = 1299.5957
$ws.Range("I68").Value = 1390
$ws.Range("J68").Value = 1288.8334
$ws.Range("K68").Value = 4170
$ws.Range("L68").Value = 3866.5002
$ws.Range("M68").Value = -3359
$ws.Range("N68").Value = -5488.5002
$ws.Range("H71").Value = 1299.5957
$ws.Range("I71").Value = 1390
$ws.Range("J71").Value = 1288.8334
$ws.Range("K71").Value = 12510
$ws.Range("L71").Value = 11599.5006
$ws.Range("M71").Value = -8454
$ws.Range("N71").Value = -19711.5006
$ws.Range("H92").Value = 340
$ws.Range("I92").Value = 350
$ws.Range("J92").Value = 333.33334
$ws.Range("K92").Value = 1050
$ws.Range("L92").Value = 1000.00002
$ws.Range("M92").Value = 198
$ws.Range("N92").Value = -3496.00002
$ws.Range("H132").Value = 2059.2
$ws.Range("I132").Value = 3230.5
$ws.Range("J132").Value = 1278.3334
$ws.Range("K132").Value = 29074.5
$ws.Range("L132").Value = 11505.0006
$ws.Range("M132").Value = -26544.5
$ws.Range("N132").Value = -16565.0006

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 4573.4
$ws.Range("I100").Value = 2380.2
$ws.Range("J100").Value = 5670
$ws.Range("K100").Value = 4760.4
$ws.Range("L100").Value = 11340
$ws.Range("M100").Value = -4219.4
$ws.Range("N100").Value = -12422
$ws.Range("H107").Value = 645.4
$ws.Range("I107").Value = 609.2222
$ws.Range("J107").Value = 699.6667
$ws.Range("K107").Value = 1827.6666
$ws.Range("L107").Value = 2099.0001
$ws.Range("M107").Value = 92.33339999999998
$ws.Range("N107").Value = -5939.0001
$ws.Range("H113").Value = 307.3125
$ws.Range("I113").Value = 278.4
$ws.Range("J113").Value = 355.5
$ws.Range("K113").Value = 835.1999999999999
$ws.Range("L113").Value = 1066.5
$ws.Range("M113").Value = 1334.8
$ws.Range("N113").Value = -5406.5
